$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.945.80"
$ws.Range("E2").Value = "  -2.71%  "

$ws.Range("D3").Value = "'2.661.01"
$ws.Range("E3").Value = "  -1.12%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'524.63"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("E6").Value = "  -1.36%  "

$ws.Range("E7").Value = "  +0.33%  "

$ws.Range("E8").Value = "  -1.13%  "

$ws.Range("E9").Value = "  +8.57%  "

$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("D13").Value = "'3.128.48"
$ws.Range("E13").Value = "  -0.88%  "

$ws.Range("D14").Value = "'58.939.84"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.671.59"
$ws.Range("E16").Value = "  -3.87%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000136"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "'338.67"
$ws.Range("E18").Value = "  -3.73%  "

$ws.Range("E19").Value = "  -3.86%  "

$ws.Range("E20").Value = "  -1.86%  "

$ws.Range("E21").Value = "  +0.65%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.06%  "

$ws.Range("D23").Value = "'63.84"

$ws.Range("E24").Value = "  -1.32%  "

$ws.Range("E25").Value = "  -2.28%  "

$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.64%  "

$ws.Range("D27").Value = "0.0₃0801"
$ws.Range("E27").Value = "  -1.99%  "

$ws.Range("D28").Value = "'7.09"
$ws.Range("E28").Value = "  -2.59%  "

$ws.Range("D29").Value = "'6.68"
$ws.Range("E29").Value = "  -3.38%  "

$ws.Range("D30").Value = "'0.999"
$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").Value = "'1.60"
$ws.Range("E31").Value = "  -0.20%  "

$ws.Range("D32").Value = "'18.85"
$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("D33").Value = "'150.52"
$ws.Range("E33").Value = "  +1.06%  "

$ws.Range("D34").Value = "'4.16"
$ws.Range("E34").Value = "  -5.24%  "

$ws.Range("E35").Value = "  -3.70%  "

$ws.Range("E36").Value = "  -6.67%  "

$ws.Range("D37").Value = "'0.870"
$ws.Range("E37").Value = "  -1.09%  "

$ws.Range("D38").Value = "'36.83"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("E39").Value = "  -6.24%  "

$ws.Range("E40").Value = "  -3.18%  "

$ws.Range("D41").Value = "'0.617"
$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("E42").Value = "  +0.24%  "

$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.95"
$ws.Range("E43").Value = "  -0.37%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").Value = "'275.38"
$ws.Range("E44").Value = "  -2.94%  "

$ws.Range("E45").Value = "  -2.38%  "

$ws.Range("E46").Value = "  +2.02%  "

$ws.Range("D47").Value = "'2.048.16"
$ws.Range("E47").Value = "  -3.77%  "

$ws.Range("E48").Value = "  -1.83%  "

$ws.Range("D49").Value = "'4.70"
$ws.Range("E49").Value = "  -3.57%  "

$ws.Range("D50").Value = "'18.92"
$ws.Range("E50").Value = "  -1.74%  "

$ws.Range("E51").Value = "  -3.15%  "
